$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (15:42 -> 16:14)
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 16:14"

# Country/provincia table refreshed with the latest COVID-19 snapshot.
# Row positions are kept fixed (sheet stays sorted by total cases, col B,
# descending); country labels and counts are rewritten in place so later
# rows correctly pick up whichever country now ranks at that position.
$rows = @(
  @{ r=4; name="Estados Unidos"; b=86548; c=1113; d=1868; e=83363; f=2122; g=22; h=1317 },
  @{ r=5; name="China"; b=81340; c=55; d=74588; e=3460; f=1034; g=5; h=3292 },
  @{ r=6; name="Italia"; b=80589; c=0; d=10361; e=62013; f=3612; g=0; h=8215 },
  @{ r=7; name="España"; b=64059; c=6273; d=9357; e=49844; f=4165; g=493; h=4858 },
  @{ r=8; name="Alemania"; b=49344; c=5406; d=5673; e=43367; f=23; g=37; h=304 },
  @{ r=9; name="Iran"; b=32332; c=2926; d=11133; e=18821; f=2893; g=144; h=2378 },
  @{ r=10; name="Francia"; b=29155; c=0; d=4948; e=22511; f=3375; g=0; h=1696 },
  @{ r=11; name="Reino Unido"; b=14579; c=2921; d=135; e=13685; f=163; g=181; h=759 },
  @{ r=12; name="Suiza"; b=12311; c=500; d=897; e=11207; f=203; g=15; h=207 },
  @{ r=13; name="Corea del Sur"; b=9332; c=91; d=4528; e=4665; f=59; g=8; h=139 },
  @{ r=14; name="Paises Bajos"; b=8603; c=1172; d=3; e=8054; f=761; g=112; h=546 },
  @{ r=15; name="Austria"; b=7441; c=532; d=225; e=7158; f=128; g=9; h=58 },
  @{ r=16; name="Belgica"; b=7284; c=1049; d=858; e=6137; f=690; g=69; h=289 },
  @{ r=17; name="Portugal"; b=4268; c=724; d=43; e=4149; f=71; g=16; h=76 },
  @{ r=18; name="Canada"; b=4043; c=0; d=228; e=3776; f=120; g=0; h=39 },
  @{ r=19; name="Noruega"; b=3694; c=322; d=6; e=3671; f=70; g=3; h=17 },
  @{ r=20; name="Turquia"; b=3629; c=0; d=26; e=3528; f=136; g=0; h=75 },
  @{ r=21; name="Australia"; b=3180; c=130; d=170; e=2997; f=23; g=0; h=13 },
  @{ r=22; name="Suecia"; b=3046; c=206; d=16; e=2938; f=214; g=15; h=92 },
  @{ r=23; name="Israel"; b=3035; c=342; d=79; e=2944; f=49; g=4; h=12 },
  @{ r=24; name="Brasil"; b=3027; c=42; d=6; e=2944; f=296; g=0; h=77 },
  @{ r=25; name="Malasia"; b=2161; c=130; d=259; e=1876; f=54; g=3; h=26 },
  @{ r=26; name="Chequia"; b=2062; c=137; d=11; e=2042; f=34; g=0; h=9 },
  @{ r=27; name="Dinamarca"; b=2010; c=133; d=1; e=1957; f=109; g=11; h=52 },
  @{ r=28; name="Irlanda"; b=1819; c=0; d=5; e=1795; f=47; g=0; h=19 },
  @{ r=29; name="Chile"; b=1610; c=304; d=43; e=1562; f=7; g=1; h=5 },
  @{ r=30; name="Luxemburgo"; b=1453; c=0; d=6; e=1438; f=3; g=0; h=9 },
  @{ r=31; name="Ecuador"; b=1403; c=0; d=3; e=1366; f=58; g=0; h=34 },
  @{ r=32; name="Japon"; b=1387; c=0; d=359; e=981; f=57; g=0; h=47 },
  @{ r=33; name="Pakistan"; b=1296; c=95; d=23; e=1264; f=7; g=0; h=9 },
  @{ r=34; name="Rumania"; b=1292; c=263; d=115; e=1153; f=32; g=1; h=24 },
  @{ r=35; name="Polonia"; b=1289; c=68; d=7; e=1266; f=3; g=0; h=16 },
  @{ r=36; name="Tailandia"; b=1136; c=91; d=97; e=1034; f=11; g=1; h=5 },
  @{ r=37; name="Arabia Saudita"; b=1104; c=92; d=35; e=1066; f=6; g=0; h=3 },
  @{ r=38; name="Indonesia"; b=1046; c=153; d=46; e=913; f=0; g=9; h=87 },
  @{ r=39; name="Finlandia"; b=1041; c=83; d=10; e=1024; f=32; g=2; h=7 },
  @{ r=40; name="Rusia"; b=1036; c=196; d=45; e=988; f=8; g=0; h=3 },
  @{ r=41; name="Sudafrica"; b=927; c=0; d=12; e=913; f=7; g=2; h=2 },
  @{ r=42; name="Grecia"; b=892; c=0; d=42; e=823; f=53; g=0; h=27 },
  @{ r=43; name="Islandia"; b=890; c=88; d=97; e=791; f=18; g=0; h=2 },
  @{ r=44; name="India"; b=863; c=136; d=73; e=770; f=0; g=0; h=20 },
  @{ r=45; name="Filipinas"; b=803; c=96; d=31; e=718; f=1; g=9; h=54 },
  @{ r=46; name="Singapur"; b=732; c=49; d=183; e=547; f=17; g=0; h=2 },
  @{ r=47; name="Crucero"; b=712; c=0; d=597; e=105; f=15; g=0; h=10 },
  @{ r=48; name="Panama"; b=674; c=0; d=2; e=663; f=20; g=0; h=9 },
  @{ r=49; name="Eslovenia"; b=632; c=70; d=10; e=613; f=14; g=3; h=9 },
  @{ r=50; name="Argentina"; b=589; c=0; d=72; e=504; f=0; g=1; h=13 },
  @{ r=51; name="Mexico"; b=585; c=110; d=4; e=573; f=1; g=2; h=8 },
  @{ r=52; name="Republica Dominicana"; b=581; c=93; d=3; e=558; f=0; g=10; h=20 },
  @{ r=53; name="Peru"; b=580; c=0; d=14; e=557; f=14; g=0; h=9 },
  @{ r=54; name="Estonia"; b=575; c=37; d=11; e=563; f=7; g=0; h=1 },
  @{ r=55; name="Croacia"; b=551; c=56; d=37; e=511; f=14; g=0; h=3 },
  @{ r=56; name="Catar"; b=549; c=0; d=43; e=506; f=6; g=0; h=0 },
  @{ r=57; name="Serbia"; b=528; c=71; d=15; e=505; f=25; g=1; h=8 },
  @{ r=58; name="Hong Kong"; b=518; c=64; d=111; e=403; f=5; g=0; h=4 },
  @{ r=59; name="Egipto"; b=495; c=0; d=102; e=369; f=0; g=0; h=24 },
  @{ r=60; name="Colombia"; b=491; c=0; d=8; e=477; f=0; g=0; h=6 },
  @{ r=61; name="Barein"; b=466; c=8; d=227; e=235; f=1; g=0; h=4 },
  @{ r=62; name="Irak"; b=458; c=76; d=122; e=296; f=0; g=4; h=40 },
  @{ r=63; name="Libano"; b=391; c=23; d=23; e=361; f=3; g=1; h=7 },
  @{ r=64; name="Nueva Zelanda"; b=368; c=0; d=37; e=331; f=1; g=0; h=0 },
  @{ r=65; name="Argelia"; b=367; c=0; d=29; e=313; f=0; g=0; h=25 },
  @{ r=66; name="Lituania"; b=345; c=46; d=1; e=339; f=2; g=1; h=5 },
  @{ r=67; name="Emiratos Arabes Unidos"; b=333; c=0; d=52; e=279; f=2; g=0; h=2 },
  @{ r=68; name="Armenia"; b=329; c=39; d=28; e=300; f=6; g=0; h=1 },
  @{ r=69; name="Hungria"; b=300; c=39; d=34; e=256; f=6; g=0; h=10 },
  @{ r=70; name="Letonia"; b=280; c=36; d=1; e=279; f=0; g=0; h=0 },
  @{ r=71; name="Bulgaria"; b=276; c=12; d=9; e=264; f=8; g=0; h=3 },
  @{ r=72; name="Marruecos"; b=275; c=0; d=8; e=256; f=1; g=0; h=11 },
  @{ r=73; name="Eslovaquia"; b=269; c=43; d=2; e=267; f=1; g=0; h=0 },
  @{ r=74; name="Principado de Andorra"; b=267; c=43; d=1; e=263; f=11; g=0; h=3 },
  @{ r=75; name="Taiwan"; b=267; c=15; d=30; e=235; f=0; g=0; h=2 },
  @{ r=76; name="Uruguay"; b=238; c=0; d=0; e=238; f=3; g=0; h=0 },
  @{ r=77; name="Costa Rica"; b=231; c=0; d=2; e=227; f=5; g=0; h=2 },
  @{ r=78; name="Bosnia y Herzegovina"; b=231; c=40; d=5; e=222; f=1; g=1; h=4 },
  @{ r=79; name="Tunez"; b=227; c=30; d=2; e=219; f=10; g=1; h=6 },
  @{ r=80; name="Ucrania"; b=226; c=30; d=5; e=216; f=0; g=0; h=5 },
  @{ r=81; name="Kuwait"; b=225; c=17; d=57; e=168; f=11; g=0; h=0 },
  @{ r=82; name="San Marino"; b=223; c=15; d=4; e=198; f=12; g=0; h=21 },
  @{ r=83; name="Republica de Macedonia"; b=219; c=18; d=3; e=213; f=1; g=0; h=3 },
  @{ r=84; name="Jordania"; b=212; c=0; d=2; e=210; f=0; g=0; h=0 },
  @{ r=85; name="Albania"; b=186; c=12; d=31; e=147; f=3; g=2; h=8 },
  @{ r=86; name="Moldavia"; b=177; c=0; d=2; e=173; f=33; g=1; h=2 },
  @{ r=87; name="Azerbaiyan"; b=165; c=43; d=15; e=147; f=6; g=0; h=3 },
  @{ r=88; name="Vietnam"; b=163; c=10; d=20; e=143; f=3; g=0; h=0 },
  @{ r=89; name="Burkina Faso"; b=152; c=0; d=10; e=135; f=0; g=0; h=7 },
  @{ r=90; name="Republica de Chipre"; b=146; c=0; d=4; e=139; f=3; g=0; h=3 },
  @{ r=91; name="Reunion"; b=145; c=10; d=1; e=144; f=0; g=0; h=0 },
  @{ r=92; name="Islas Feroe"; b=144; c=4; d=54; e=90; f=2; g=0; h=0 },
  @{ r=93; name="Malta"; b=139; c=5; d=2; e=137; f=1; g=0; h=0 },
  @{ r=94; name="Ghana"; b=136; c=4; d=1; e=131; f=1; g=0; h=4 },
  @{ r=95; name="Kazajistan"; b=135; c=22; d=3; e=131; f=0; g=0; h=1 },
  @{ r=96; name="Oman"; b=131; c=22; d=23; e=108; f=0; g=0; h=0 },
  @{ r=97; name="Senegal"; b=119; c=14; d=11; e=108; f=0; g=0; h=0 },
  @{ r=98; name="Brunei"; b=115; c=1; d=11; e=104; f=1; g=0; h=0 },
  @{ r=99; name="Venezuela"; b=107; c=0; d=31; e=75; f=2; g=0; h=1 },
  @{ r=100; name="Sri Lanka"; b=106; c=0; d=7; e=99; f=5; g=0; h=0 },
  @{ r=101; name="Camboya"; b=99; c=1; d=11; e=88; f=1; g=0; h=0 },
  @{ r=102; name="Costa de Marfil"; b=96; c=0; d=3; e=93; f=0; g=0; h=0 },
  @{ r=103; name="Mauricio"; b=94; c=13; d=0; e=92; f=1; g=0; h=2 },
  @{ r=104; name="Afganistan"; b=94; c=0; d=2; e=88; f=0; g=0; h=4 },
  @{ r=105; name="Bielorrusia"; b=94; c=8; d=32; e=62; f=2; g=0; h=0 },
  @{ r=106; name="Estado de Palestina"; b=91; c=5; d=17; e=73; f=0; g=0; h=1 },
  @{ r=107; name="Camerun"; b=88; c=13; d=2; e=84; f=0; g=1; h=2 },
  @{ r=108; name="Uzbekistan"; b=83; c=8; d=5; e=77; f=8; g=1; h=1 },
  @{ r=109; name="Martinica"; b=81; c=0; d=0; e=80; f=12; g=0; h=1 },
  @{ r=110; name="Georgia"; b=81; c=2; d=13; e=68; f=1; g=0; h=0 },
  @{ r=111; name="Guadalupe"; b=73; c=0; d=0; e=72; f=4; g=0; h=1 },
  @{ r=112; name="Montenegro"; b=70; c=1; d=0; e=69; f=1; g=0; h=1 },
  @{ r=113; name="Honduras"; b=68; c=1; d=0; e=67; f=0; g=0; h=1 },
  @{ r=114; name="Cuba"; b=67; c=0; d=1; e=64; f=2; g=0; h=2 },
  @{ r=115; name="Trinidad yTobago"; b=66; c=1; d=1; e=63; f=0; g=1; h=2 },
  @{ r=116; name="Nigeria"; b=65; c=0; d=3; e=61; f=0; g=0; h=1 },
  @{ r=117; name="Bolivia"; b=61; c=0; d=0; e=61; f=0; g=0; h=0 },
  @{ r=118; name="Kirguistan"; b=58; c=14; d=0; e=58; f=0; g=0; h=0 },
  @{ r=119; name="Liechtenstein"; b=56; c=0; d=0; e=56; f=0; g=0; h=0 },
  @{ r=120; name="Paraguay"; b=52; c=11; d=1; e=48; f=1; g=0; h=3 },
  @{ r=121; name="Consejo Danes para los Refugiados"; b=51; c=0; d=2; e=46; f=0; g=0; h=3 },
  @{ r=122; name="Mayotte"; b=50; c=14; d=0; e=50; f=0; g=0; h=0 },
  @{ r=123; name="Ruanda"; b=50; c=0; d=0; e=50; f=0; g=0; h=0 },
  @{ r=124; name="Banglades"; b=48; c=4; d=11; e=32; f=1; g=0; h=5 },
  @{ r=125; name="Puerto Rico"; b=39; c=0; d=1; e=36; f=0; g=0; h=2 },
  @{ r=126; name="Gibraltar"; b=35; c=0; d=13; e=22; f=0; g=0; h=0 },
  @{ r=127; name="Macao"; b=34; c=1; d=10; e=24; f=0; g=0; h=0 },
  @{ r=128; name="Monaco"; b=33; c=0; d=1; e=32; f=0; g=0; h=0 },
  @{ r=129; name="Guam"; b=32; c=0; d=0; e=31; f=0; g=0; h=1 },
  @{ r=130; name="Kenia"; b=31; c=0; d=1; e=29; f=0; g=0; h=1 },
  @{ r=131; name="Polinesia Francesa"; b=30; c=0; d=0; e=30; f=0; g=0; h=0 },
  @{ r=132; name="Isla de Man"; b=29; c=3; d=0; e=29; f=0; g=0; h=0 },
  @{ r=133; name="Aruba"; b=28; c=0; d=1; e=27; f=0; g=0; h=0 },
  @{ r=134; name="Guayana Francesa"; b=28; c=0; d=6; e=22; f=0; g=0; h=0 },
  @{ r=135; name="Jamaica"; b=26; c=0; d=2; e=23; f=0; g=0; h=1 },
  @{ r=136; name="Togo"; b=25; c=1; d=1; e=24; f=0; g=0; h=0 },
  @{ r=137; name="Guatemala"; b=25; c=0; d=4; e=20; f=0; g=0; h=1 },
  @{ r=138; name="Barbados"; b=24; c=0; d=0; e=24; f=0; g=0; h=0 },
  @{ r=139; name="Madagascar"; b=24; c=1; d=0; e=24; f=0; g=0; h=0 },
  @{ r=140; name="Zambia"; b=22; c=6; d=0; e=22; f=0; g=0; h=0 },
  @{ r=141; name="Uganda"; b=18; c=4; d=0; e=18; f=0; g=0; h=0 },
  @{ r=142; name="Islas Virgenes de los Estados Unidos"; b=17; c=0; d=0; e=17; f=0; g=0; h=0 },
  @{ r=143; name="Etiopia"; b=16; c=4; d=0; e=16; f=0; g=0; h=0 },
  @{ r=144; name="Nueva Caledonia"; b=15; c=1; d=0; e=15; f=0; g=0; h=0 },
  @{ r=145; name="Bermudas"; b=15; c=0; d=2; e=13; f=0; g=0; h=0 },
  @{ r=146; name="Maldivas"; b=14; c=1; d=9; e=5; f=0; g=0; h=0 },
  @{ r=147; name="El Salvador"; b=13; c=0; d=0; e=13; f=0; g=0; h=0 },
  @{ r=148; name="Tanzania"; b=13; c=0; d=1; e=12; f=0; g=0; h=0 },
  @{ r=149; name="Guinea Ecuatorial"; b=12; c=0; d=0; e=12; f=0; g=0; h=0 },
  @{ r=150; name="Republica de Yibuti"; b=12; c=1; d=0; e=12; f=0; g=0; h=0 },
  @{ r=151; name="San Martin (Parte Francesa)"; b=11; c=0; d=0; e=11; f=0; g=0; h=0 },
  @{ r=152; name="Mongolia"; b=11; c=0; d=0; e=11; f=0; g=0; h=0 },
  @{ r=153; name="Dominica"; b=11; c=0; d=0; e=11; f=0; g=0; h=0 },
  @{ r=154; name="Mali"; b=11; c=7; d=0; e=11; f=0; g=0; h=0 },
  @{ r=155; name="Niger"; b=10; c=0; d=0; e=9; f=0; g=0; h=1 },
  @{ r=156; name="Bahamas"; b=9; c=0; d=1; e=8; f=0; g=0; h=0 },
  @{ r=157; name="Groenlandia"; b=9; c=3; d=2; e=7; f=0; g=0; h=0 },
  @{ r=158; name="Surinam"; b=8; c=0; d=0; e=8; f=0; g=0; h=0 },
  @{ r=159; name="Guinea"; b=8; c=4; d=0; e=8; f=0; g=0; h=0 },
  @{ r=160; name="Haiti"; b=8; c=0; d=0; e=8; f=0; g=0; h=0 },
  @{ r=161; name="Islas Caimanes"; b=8; c=0; d=0; e=7; f=0; g=0; h=1 },
  @{ r=162; name="Namibia"; b=8; c=0; d=2; e=6; f=0; g=0; h=0 },
  @{ r=163; name="Mozambique"; b=7; c=0; d=0; e=7; f=0; g=0; h=0 },
  @{ r=164; name="Seychelles"; b=7; c=0; d=0; e=7; f=0; g=0; h=0 },
  @{ r=165; name="Antigua y Barbuda"; b=7; c=0; d=0; e=7; f=0; g=0; h=0 },
  @{ r=166; name="Granada"; b=7; c=0; d=0; e=7; f=0; g=0; h=0 },
  @{ r=167; name="Gabon"; b=7; c=0; d=0; e=6; f=0; g=0; h=1 },
  @{ r=168; name="Curazao"; b=7; c=0; d=2; e=4; f=0; g=0; h=1 },
  @{ r=169; name="Suazilandia"; b=6; c=0; d=0; e=6; f=0; g=0; h=0 },
  @{ r=170; name="Benin"; b=6; c=0; d=0; e=6; f=0; g=0; h=0 },
  @{ r=171; name="Laos"; b=6; c=0; d=0; e=6; f=0; g=0; h=0 },
  @{ r=172; name="Eritrea"; b=6; c=0; d=0; e=6; f=0; g=0; h=0 },
  @{ r=173; name="Siria"; b=5; c=0; d=0; e=5; f=0; g=0; h=0 },
  @{ r=174; name="Birmania"; b=5; c=0; d=0; e=5; f=0; g=0; h=0 },
  @{ r=175; name="Fiyi"; b=5; c=0; d=0; e=5; f=0; g=0; h=0 },
  @{ r=176; name="Montserrat"; b=5; c=0; d=0; e=5; f=0; g=0; h=0 },
  @{ r=177; name="Cabo Verde"; b=5; c=0; d=0; e=4; f=0; g=0; h=1 },
  @{ r=178; name="Guyana"; b=5; c=0; d=0; e=4; f=0; g=0; h=1 },
  @{ r=179; name="Zimbabue"; b=5; c=2; d=0; e=4; f=0; g=0; h=1 },
  @{ r=180; name="Congo"; b=4; c=0; d=0; e=4; f=0; g=0; h=0 },
  @{ r=181; name="Santa Sede"; b=4; c=0; d=0; e=4; f=0; g=0; h=0 },
  @{ r=182; name="Angola"; b=4; c=0; d=0; e=4; f=0; g=0; h=0 },
  @{ r=183; name="Nepal"; b=4; c=1; d=1; e=3; f=0; g=0; h=0 },
  @{ r=184; name="Mauritania"; b=3; c=0; d=0; e=3; f=0; g=0; h=0 },
  @{ r=185; name="Republica de Africa Central"; b=3; c=0; d=0; e=3; f=0; g=0; h=0 },
  @{ r=186; name="Liberia"; b=3; c=0; d=0; e=3; f=0; g=0; h=0 },
  @{ r=187; name="San Martin (Parte Holandesa)"; b=3; c=0; d=0; e=3; f=0; g=0; h=0 },
  @{ r=188; name="Republica del Chad"; b=3; c=0; d=0; e=3; f=0; g=0; h=0 },
  @{ r=189; name="Butan"; b=3; c=1; d=0; e=3; f=0; g=0; h=0 },
  @{ r=190; name="San Bartolome"; b=3; c=0; d=0; e=3; f=0; g=0; h=0 },
  @{ r=191; name="Somalia"; b=3; c=1; d=0; e=3; f=0; g=0; h=0 },
  @{ r=192; name="Gambia"; b=3; c=0; d=0; e=2; f=0; g=0; h=1 },
  @{ r=193; name="Sudan"; b=3; c=0; d=0; e=2; f=0; g=0; h=1 },
  @{ r=194; name="Santa Lucia"; b=3; c=0; d=1; e=2; f=0; g=0; h=0 },
  @{ r=195; name="Islas Virgenes Britanicas"; b=2; c=0; d=0; e=2; f=0; g=0; h=0 },
  @{ r=196; name="San Cristobal y Nieves"; b=2; c=0; d=0; e=2; f=0; g=0; h=0 },
  @{ r=197; name="Guinea-Bisau"; b=2; c=0; d=0; e=2; f=0; g=0; h=0 },
  @{ r=198; name="Anguila"; b=2; c=0; d=0; e=2; f=0; g=0; h=0 },
  @{ r=199; name="Islas Turcas y Caicos"; b=2; c=0; d=0; e=2; f=0; g=0; h=0 },
  @{ r=200; name="Belice"; b=2; c=0; d=0; e=2; f=0; g=0; h=0 },
  @{ r=201; name="Nicaragua"; b=2; c=0; d=0; e=1; f=0; g=1; h=1 },
  @{ r=202; name="San Vicente y las Granadinas"; b=1; c=0; d=0; e=1; f=0; g=0; h=0 },
  @{ r=203; name="Libia"; b=1; c=0; d=0; e=1; f=0; g=0; h=0 },
  @{ r=204; name="Timor Oriental"; b=1; c=0; d=0; e=1; f=0; g=0; h=0 },
  @{ r=205; name="Papua Nueva Guinea"; b=1; c=0; d=0; e=1; f=0; g=0; h=0 }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value = $row.name
    $ws.Cells.Item($row.r, 2).Value = $row.b
    $ws.Cells.Item($row.r, 3).Value = $row.c
    $ws.Cells.Item($row.r, 4).Value = $row.d
    $ws.Cells.Item($row.r, 5).Value = $row.e
    $ws.Cells.Item($row.r, 6).Value = $row.f
    $ws.Cells.Item($row.r, 7).Value = $row.g
    $ws.Cells.Item($row.r, 8).Value = $row.h
}
